$d = $word.ActiveDocument

# wdAlignParagraphLeft = 0, wdAlignParagraphJustify = 3
# Remove the "justify" (both-sides) alignment from every paragraph that
# currently has it, restoring the default (unset) paragraph alignment -
# i.e. drop <w:jc w:val="both"/> from <w:pPr>.
foreach ($p in $d.Paragraphs) {
    if ($p.Format.Alignment -eq 3) {
        $p.Format.Alignment = 0
    }
}
